$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "<are>"
$ws.Range("B2").Value = "<are>"
$ws.Range("C2").Value = 6

$ws.Range("A3").Value = "<his>"
$ws.Range("B3").Value = "<is>"
$ws.Range("C3").Value = 14

$ws.Range("A4").Value = "<kilo>"
$ws.Range("B4").Value = "<kilo>"
$ws.Range("C4").Value = 4

$ws.Range("A5").Value = "<come>"
$ws.Range("B5").Value = "<come>"
$ws.Range("C5").Value = 8

$ws.Range("A6").Value = "<sentence>"
$ws.Range("B6").Value = "<sentence>"
$ws.Range("C6").Value = 11

$ws.Range("A7").Value = "<zero>"
$ws.Range("B7").Value = "<zero>"
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = "<shift>"
$ws.Range("B8").Value = "<shift>"
$ws.Range("C8").Value = 7

$ws.Range("A9").Value = "<so>"
$ws.Range("B9").Value = "<so>"
$ws.Range("C9").Value = 4

$ws.Range("A10").Value = "<lima>"
$ws.Range("B10").Value = "<lima>"
$ws.Range("C10").Value = 9

$ws.Range("A11").Value = "<be>"
$ws.Range("B11").Value = "<be>"
$ws.Range("C11").Value = 8

$ws.Range("A12").Value = "<him>"
$ws.Range("B12").Value = "<him>"
$ws.Range("C12").Value = 1

$ws.Range("A13").Value = "<your>"
$ws.Range("B13").Value = "<your>"
$ws.Range("C13").Value = 9

$ws.Range("A14").Value = "<are>"
$ws.Range("B14").Value = "<are>"
$ws.Range("C14").Value = 7

$ws.Range("A15").Value = "<in>"
$ws.Range("B15").Value = "<in>"
$ws.Range("C15").Value = 12

$ws.Range("A16").Value = "<number>"
$ws.Range("B16").Value = "<number>"
$ws.Range("C16").Value = 8

$ws.Range("A17").Value = "<left>"
$ws.Range("B17").Value = "<left>"
$ws.Range("C17").Value = 8

$ws.Range("A18").Value = "<no>"
$ws.Range("B18").Value = "<no>"
$ws.Range("C18").Value = 7
